$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "25/03/2023"

$data = @{
  2  = @(1, 844, 974, 32, 3, 14, 94, 1666, -41.53661464585834)
  3  = @(1, 75, 77, 2, 0, 3, 0, 107, -28.0373831775701)
  4  = @(0, 182, 183, 2, 0, 3, 0, 63, 190.4761904761905)
  5  = @(0, 672, 747, 30, 4, 11, 41, 821, -9.013398294762487)
  6  = @(1, 369, 433, 70, 1, 7, 0, 564, -23.22695035460993)
  7  = @(0, 180, 190, 7, 3, 2, 0, 174, 9.195402298850585)
  8  = @(0, 63, 92, 29, 0, 0, 0, 120, -23.33333333333333)
  9  = @(0, 323, 335, 9, 2, 3, 1, 778, -56.94087403598972)
  10 = @(0, 39, 37, 0, 0, 0, 0, 73, -49.31506849315068)
  11 = @(1, 0, 1, 0, 0, 0, 0, 3, -66.66666666666667)
  12 = @(0, 37, 37, 0, 0, 2, 0, 52, -28.84615384615384)
}

foreach ($row in $data.Keys) {
  $values = $data[$row]
  for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item($row, $col).Value = $values[$i]
  }
}
